$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2130.0908
$ws.Range("J40").Value = 2394.75
$ws.Range("L40").Value = 2394.75
$ws.Range("N40").Value = -2744.75
$ws.Range("H55").Value = 184.66667
$ws.Range("J55").Value = 184
$ws.Range("L55").Value = 184
$ws.Range("N55").Value = -612
$ws.Range("H98").Value = 1898.0322
$ws.Range("I98").Value = 1605.7693
$ws.Range("K98").Value = 1605.7693
$ws.Range("M98").Value = -107.7692999999999
$ws.Range("H112").Value = 1817.7759
$ws.Range("I112").Value = 1500
$ws.Range("K112").Value = 4500
$ws.Range("M112").Value = -3392
$ws.Range("H122").Value = 1898.0322
$ws.Range("I122").Value = 1605.7693
$ws.Range("K122").Value = 4817.3079
$ws.Range("M122").Value = -2367.3079
$ws.Range("H129").Value = 1220.3448
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1220.3448
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3661.0344
$ws.Range("M129").Value = $null
$ws.Range("N129").Value = -13661.0344
$ws.Range("H132").Value = 1702.4
$ws.Range("I132").Value = 1900
$ws.Range("J132").Value = 1570.6666
$ws.Range("K132").Value = 5700
$ws.Range("L132").Value = 4711.9998
$ws.Range("M132").Value = -3170
$ws.Range("N132").Value = -9771.9998
$ws.Range("H138").Value = 2646.4211
$ws.Range("J138").Value = 2428.7856
$ws.Range("L138").Value = 7286.3568
$ws.Range("N138").Value = -17566.3568
$ws.Range("H139").Value = 73600
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").Value = $null
$ws.Range("H140").Value = 81792.234
$ws.Range("J140").Value = 81792.234
$ws.Range("L140").Value = 81792.234
$ws.Range("N140").Value = -92152.234
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 463479
$ws.Range("I2").Value = 694954.5
$ws.Range("K2").Value = 694954.5
$ws.Range("M2").Value = -694841.5
$ws.Range("H32").Value = 3708.2856
$ws.Range("I32").Value = 3062.3152
$ws.Range("J32").Value = 15497.25
$ws.Range("K32").Value = 3062.3152
$ws.Range("L32").Value = 15497.25
$ws.Range("M32").Value = -2775.3152
$ws.Range("N32").Value = -16071.25
$ws.Range("H116").Value = 463479
$ws.Range("I116").Value = 694954.5
$ws.Range("K116").Value = 694954.5
$ws.Range("M116").Value = -692660.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 463479
$ws.Range("I3").Value = 694954.5
$ws.Range("K3").Value = 694954.5
$ws.Range("M3").Value = -694840.5
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
$ws.Range("H134").Value = 5666.926
$ws.Range("I134").Value = 5975.2
$ws.Range("K134").Value = 17925.6
$ws.Range("M134").Value = -15390.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3438.818
$ws.Range("I31").Value = 3812.1667
$ws.Range("K31").Value = 3812.1667
$ws.Range("M31").Value = -3517.1667
$ws.Range("H34").Value = 3438.818
$ws.Range("I34").Value = 3812.1667
$ws.Range("K34").Value = 3812.1667
$ws.Range("M34").Value = -3610.1667
$ws.Range("H141").Value = 55674.75
$ws.Range("J141").Value = 53342.57
$ws.Range("L141").Value = 53342.57
$ws.Range("N141").Value = -63702.57
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 644.64703
$ws.Range("I5").Value = 557
$ws.Range("J5").Value = 855
$ws.Range("K5").Value = 1671
$ws.Range("L5").Value = 2565
$ws.Range("M5").Value = -1559
$ws.Range("N5").Value = -2789
$ws.Range("H44").Value = 2099.5
$ws.Range("I44").Value = 199
$ws.Range("J44").Value = 4000
$ws.Range("K44").Value = 597
$ws.Range("L44").Value = 12000
$ws.Range("M44").Value = -199
$ws.Range("N44").Value = -12796
$ws.Range("H68").Value = 167286
$ws.Range("I68").Value = 400
$ws.Range("J68").Value = 200663.2
$ws.Range("K68").Value = 1200
$ws.Range("L68").Value = 601989.6000000001
$ws.Range("M68").Value = -389
$ws.Range("N68").Value = -603611.6000000001
$ws.Range("H71").Value = 167286
$ws.Range("I71").Value = 400
$ws.Range("J71").Value = 200663.2
$ws.Range("K71").Value = 3600
$ws.Range("L71").Value = 1805968.8
$ws.Range("M71").Value = 456
$ws.Range("N71").Value = -1814080.8
$ws.Range("H81").Value = 1515.5714
$ws.Range("I81").Value = 920
$ws.Range("K81").Value = 2760
$ws.Range("M81").Value = -1637
$ws.Range("H84").Value = 1515.5714
$ws.Range("I84").Value = 920
$ws.Range("K84").Value = 8280
$ws.Range("M84").Value = -2664
$ws.Range("H98").Value = 200976.6
$ws.Range("J98").Value = 200976.6
$ws.Range("L98").Value = 602929.8
$ws.Range("N98").Value = -605925.8
$ws.Range("H104").Value = 4500.077
$ws.Range("J104").Value = 4500.077
$ws.Range("L104").Value = 13500.231
$ws.Range("N104").Value = -18742.231
$ws.Range("H107").Value = 440.63635
$ws.Range("J107").Value = 405.22223
$ws.Range("L107").Value = 1215.66669
$ws.Range("N107").Value = -5055.66669
$ws.Range("H113").Value = 5508.905
$ws.Range("J113").Value = 732.55554
$ws.Range("L113").Value = 2197.66662
$ws.Range("N113").Value = -6537.66662
$ws.Range("H125").Value = 4885.7144
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 5533.3335
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 16600.0005
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -26440.0005
$ws.Range("H131").Value = 17798.781
$ws.Range("J131").Value = 19150.525
$ws.Range("L131").Value = 57451.575
$ws.Range("N131").Value = -67531.57500000001
$ws.Range("H135").Value = 644.64703
$ws.Range("I135").Value = 557
$ws.Range("J135").Value = 855
$ws.Range("K135").Value = 5013
$ws.Range("L135").Value = 7695
$ws.Range("M135").Value = -2478
$ws.Range("N135").Value = -12765
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3537965
$ws.Range("J126").Value = 252618.5
$ws.Range("L126").Value = 757855.5
$ws.Range("N126").Value = -762795.5
$ws.Range("H132").Value = 3848715
$ws.Range("I132").Value = 6411758.5
$ws.Range("J132").Value = 4149.75
$ws.Range("K132").Value = 19235275.5
$ws.Range("L132").Value = 12449.25
$ws.Range("M132").Value = -19232745.5
$ws.Range("N132").Value = -17509.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = $null
$ws.Range("H122").Value = 38953.145
$ws.Range("I122").Value = 88043.11
$ws.Range("K122").Value = 264129.33
$ws.Range("M122").Value = -261679.33
$ws.Range("H123").Value = 61746.547
$ws.Range("J123").Value = 61746.547
$ws.Range("L123").Value = 61746.547
$ws.Range("N123").Value = -71546.54699999999
$ws.Range("H133").Value = 63873.75
$ws.Range("I133").Value = 63715
$ws.Range("J133").Value = 63926.668
$ws.Range("K133").Value = 63715
$ws.Range("L133").Value = 63926.668
$ws.Range("M133").Value = -58655
$ws.Range("N133").Value = -74046.66800000001
$ws.Range("H136").Value = 29241590
$ws.Range("I136").Value = 61729810
$ws.Range("J136").Value = 2194
$ws.Range("K136").Value = 185189430
$ws.Range("L136").Value = 6582
$ws.Range("M136").Value = -185186880
$ws.Range("N136").Value = -11682
